$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Foosball")
$ws.Range("H10").Value = 7
$ws.Range("H11").Value = 7
$ws.Range("H12").Value = 5
$ws.Range("H13").Value = 2
$ws.Range("H14").Value = 2
$ws.Range("H17").Value = 4
$ws.Range("H18").Value = 3
$ws.Range("H19").Value = 5
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 6
$ws.Range("H22").Value = 3
$ws.Range("H23").Value = 1
$ws.Range("H24").Value = 4
$ws.Range("H25").Value = 6

$ws = $wb.Worksheets.Item("Carrom")
$ws.Range("H10").Value = 7
$ws.Range("H11").Value = 5
$ws.Range("H12").Value = 6
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 8
$ws.Range("H15").Value = 8
$ws.Range("H16").Value = 4
$ws.Range("H17").Value = 2
$ws.Range("H18").Value = 4
$ws.Range("H19").Value = 3
$ws.Range("H20").Value = 7
$ws.Range("H21").Value = 2
$ws.Range("H22").Value = 1
$ws.Range("H23").Value = 6
$ws.Range("H24").Value = 5
$ws.Range("H25").Value = 3

$ws = $wb.Worksheets.Item("Badminton")
$ws.Range("H10").Value = 3
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 8
$ws.Range("H13").Value = 5
$ws.Range("H14").Value = 4
$ws.Range("H15").Value = 7
$ws.Range("H16").Value = 1
$ws.Range("H17").Value = 6
$ws.Range("H18").Value = 2
$ws.Range("H19").Value = 6
$ws.Range("H20").Value = 4
$ws.Range("H21").Value = 5
$ws.Range("H22").Value = 3
$ws.Range("H23").Value = 2
$ws.Range("H24").Value = 8
$ws.Range("H25").Value = 7

$ws = $wb.Worksheets.Item("Chess")
$ws.Range("G18").Value = 15
$ws.Range("G21").Value = 9
$ws.Range("G22").Value = 3
$ws.Range("G23").Value = 2
$ws.Range("G24").Value = 13
$ws.Range("G25").Value = 7
$ws.Range("G26").Value = 13
$ws.Range("G27").Value = 11
$ws.Range("G28").Value = 8
$ws.Range("G29").Value = 14
$ws.Range("G30").Value = 6
$ws.Range("G31").Value = 10
$ws.Range("G32").Value = 12
$ws.Range("G33").Value = 5
$ws.Range("G34").Value = 1
$ws.Range("G35").Value = 12
$ws.Range("G36").Value = 3
$ws.Range("G37").Value = 9
$ws.Range("G38").Value = 6
$ws.Range("G39").Value = 4
$ws.Range("G40").Value = 10
$ws.Range("G41").Value = 14
$ws.Range("G42").Value = 11
$ws.Range("G43").Value = 4
$ws.Range("G44").Value = 2
$ws.Range("G45").Value = 7
$ws.Range("G46").Value = 15
$ws.Range("G47").Value = 5
$ws.Range("G48").Value = 8
$ws.Range("G49").Value = 1

$ws = $wb.Worksheets.Item("Table tennis")
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 8
$ws.Range("H12").Value = 5
$ws.Range("H13").Value = 3
$ws.Range("H14").Value = 8
$ws.Range("H15").Value = 6
$ws.Range("H16").Value = 7
$ws.Range("H17").Value = 4
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 3
$ws.Range("H20").Value = 5
$ws.Range("H21").Value = 2
$ws.Range("H22").Value = 4
$ws.Range("H23").Value = 7
$ws.Range("H24").Value = 2
$ws.Range("H25").Value = 6

$wb.Worksheets.Item("Foosball").Activate()
